$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (Förändrad) rows 2-308: update date serial value from 45206 to 45208
$ws.Range("C2:C308").Value = 45208
